$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update "Förändrad" (changed) date in column C for rows 2-8
# from serial 45185 (2023-09-16) to serial 45204 (2023-10-05)
for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
